$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet holds a weekly price history for "Betarraga" (Vega Monumental
# Concepcion). A new week's observation (date 44476) is prepended to the
# existing Primera/Segunda row-pairs (rows 188-207), pushing every older
# pair down by one slot; the oldest pair that falls off the bottom of the
# block is appended as two brand-new rows (208-209).
#
# Only columns D (Fecha) and J (Volumen) differ between the shifted pairs;
# every other column is constant for a given quality (Primera/Segunda), so
# we update D/J in place for rows 188-207 and then write full rows for the
# two newly appended rows 208-209.

$dates = @{
    188 = 44476; 189 = 44476
    190 = 44386; 191 = 44386
    192 = 44306; 193 = 44306
    194 = 44357; 195 = 44357
    196 = 44321; 197 = 44321
    198 = 44397; 199 = 44397
    200 = 44314; 201 = 44314
    202 = 44425; 203 = 44425
    204 = 44390; 205 = 44390
    206 = 44250; 207 = 44250
}

$volumes = @{
    188 = 100; 189 = 50
    190 = 800; 191 = 400
    192 = 600; 193 = 300
    194 = 800; 195 = 400
    196 = 600; 197 = 300
    198 = 600; 199 = 300
    200 = 400; 201 = 200
    202 = 600; 203 = 300
    204 = 600; 205 = 300
    206 = 600; 207 = 300
}

foreach ($r in 188..207) {
    $ws.Cells.Item($r, 4).Value = $dates[$r]
    $ws.Cells.Item($r, 10).Value = $volumes[$r]
}

function Set-PriceRow($row, $quality, $volumen, $kmin, $kmax, $kprom, $precioKg) {
    $ws.Cells.Item($row, 1).Value = 11
    $ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
    $ws.Cells.Item($row, 3).Value = "Bíobío"
    $ws.Cells.Item($row, 4).Value = 44432
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = 8
    $ws.Cells.Item($row, 6).Value = 100114014
    $ws.Cells.Item($row, 7).Value = "Betarraga"
    $ws.Cells.Item($row, 8).Value = "Sin especificar"
    $ws.Cells.Item($row, 9).Value = $quality
    $ws.Cells.Item($row, 10).Value = $volumen
    $ws.Cells.Item($row, 11).Value = $kmin
    $ws.Cells.Item($row, 12).Value = $kmax
    $ws.Cells.Item($row, 13).Value = $kprom
    $ws.Cells.Item($row, 14).Value = "`$/paquete 5 unidades"
    $ws.Cells.Item($row, 15).Value = "Región Metropolitana"
    $ws.Cells.Item($row, 16).Value = $precioKg
    $ws.Cells.Item($row, 17).Value = 5
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}

# New rows appended at the bottom, carrying the pair that fell off the
# shifted block (old row 206/207 values: date 44432).
Set-PriceRow 208 "Primera" 600 600 700 650 130
Set-PriceRow 209 "Segunda" 300 500 500 500 100
